$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 85, shifting existing rows 85-89 down to 86-90.
$ws.Rows.Item(85).Insert()

# Match the date-column style used by the rest of the table (row now at 86,
# formerly row 85) for the newly inserted row's date cell.
$ws.Cells.Item(86, 4).Copy()
$ws.Cells.Item(85, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(85, 1).Value = 3
$ws.Cells.Item(85, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44826
$ws.Cells.Item(85, 5).Value = 5
$ws.Cells.Item(85, 6).Value = 100112035
$ws.Cells.Item(85, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 50
$ws.Cells.Item(85, 11).Value = 15000
$ws.Cells.Item(85, 12).Value = 15000
$ws.Cells.Item(85, 13).Value = 15000
$ws.Cells.Item(85, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(85, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(85, 16).Value = 1000
$ws.Cells.Item(85, 17).Value = 15
$ws.Cells.Item(85, 18).Value = "Hortaliza"
